$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1682
$ws1.Range("F5").Value = 2276
$ws1.Range("F9").Value = 1042
$ws1.Range("F10").Value = 259
$ws1.Range("F16").Value = 8059
$ws1.Range("F17").Value = 352
$ws1.Range("F19").Value = 225
$ws1.Range("F28").Value = 437
$ws1.Range("F35").Value = 55

# Sheet "全部类型" (all types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1682
$ws4.Range("F7").Value = 2276
$ws4.Range("F12").Value = 1042
$ws4.Range("F13").Value = 259
$ws4.Range("F18").Value = 8059
$ws4.Range("F19").Value = 352
$ws4.Range("F22").Value = 225
$ws4.Range("F31").Value = 437
$ws4.Range("F38").Value = 55
